$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Full_Data")

# Delete row 74 (PBS / SRR796591 entry) - everything below shifts up by one.
$ws.Rows(74).Delete()
